# Scene.xlsx server-data update
# - Removes the "CloneScene" / "Scene2" row
# - Removes the "RebellerNoob" / "SelectScene" row
# - Updates the RelivePos for the PioneerNoob/villageScene row to "20,0,60"
# - Updates the ID for the Demo1 row to "2"
# - Leaves the selection on F5, matching the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "CloneScene" row (originally row 2: FilePath=".../CloneScene/", SceneName="Scene2")
$ws.Rows(2).Delete()

# After the previous delete, the "RebellerNoob" row (FilePath=".../RebellerNoob/",
# SceneName="SelectScene") has shifted up to row 3 - remove it as well.
$ws.Rows(3).Delete()

# Remaining row 2 is now PioneerNoob / villageScene - update its RelivePos value.
$ws.Range("E2").Value2 = "20,0,60"

# Remaining row 3 is now the Demo1 row - update its ID value.
$ws.Range("B3").Value2 = "2"

# Match the saved selection state from the workbook.
$ws.Range("F5").Select() | Out-Null
